$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - Regular US Data
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 126
$ws.Range("E6").Value = 3331000360
$ws.Range("F6").Value = 5
$ws.Range("G6").Formula = "=(`$B`$2/B6+`$B`$2/D6-`$B`$2/C6)*F6"
$ws.Range("H6").Formula = "=G6+E6"
$ws.Range("I6").Formula = "=IF(H6>E6, ""Increase"", ""Decrease"")"

# Row 7 - Population Increase - High Birth Rate
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = 126
$ws.Range("E7").Value = 3331000360
$ws.Range("F7").Value = 5
$ws.Range("G7").Formula = "=(`$B`$2/B7+`$B`$2/D7-`$B`$2/C7)*F7"
$ws.Range("H7").Formula = "=G7+E7"
$ws.Range("I7").Formula = "=IF(H7>E7, ""Increase"", ""Decrease"")"

# Row 8 - Population Increase - High Migration
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 12
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 3331000360
$ws.Range("F8").Value = 5
$ws.Range("G8").Formula = "=(`$B`$2/B8+`$B`$2/D8-`$B`$2/C8)*F8"
$ws.Range("H8").Formula = "=G8+E8"
$ws.Range("I8").Formula = "=IF(H8>E8, ""Increase"", ""Decrease"")"

# Row 9 - Population Decrease - High Death Rate
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 126
$ws.Range("E9").Value = 3331000360
$ws.Range("F9").Value = 5
$ws.Range("G9").Formula = "=(`$B`$2/B9+`$B`$2/D9-`$B`$2/C9)*F9"
$ws.Range("H9").Formula = "=G9+E9"
$ws.Range("I9").Formula = "=IF(H9>E9, ""Increase"", ""Decrease"")"

# Row 10 - Population Low Birth Rate and Low Migration
$ws.Range("B10").Value = 50
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 200
$ws.Range("E10").Value = 3331000360
$ws.Range("F10").Value = 5
$ws.Range("G10").Formula = "=(`$B`$2/B10+`$B`$2/D10-`$B`$2/C10)*F10"
$ws.Range("H10").Formula = "=G10+E10"
$ws.Range("I10").Formula = "=IF(H10>E10, ""Increase"", ""Decrease"")"

# Update selection to match final state
$ws.Range("H15").Select()
